$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("potentiometers")

# Update the formula in column A (rows 2:22) to use 4-digit zero padding
# instead of 6-digit, e.g. "POT-000001" -> "POT-0001"
$ws.Range("A2").Formula = '="POT-"&TEXT(ROW()-1,"0000")'
$ws.Range("A3:A22").Formula = '="POT-"&TEXT(ROW()-1,"0000")'

# Update the active selection on the sheet: previously K2, now A2:A22
$ws.Range("A2:A22").Select()
